$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.888.65"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "3.788.50"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").Value = "'602.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "'163.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").Value = "3.787.64"
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("D11").Value = "'6.94"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +10.11%  "
$ws.Range("D12").Value = "'0.446"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").Value = "'35.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "4.422.71"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "3.792.52"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "67.883.27"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "'18.19"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").Value = "'7.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").Value = "'458.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").Value = "'9.44"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.07%  "
$ws.Range("D23").Value = "'0.691"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").Value = "'83.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  -4.61%  "
$ws.Range("D26").Value = "'11.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "'9.93"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").Value = "3.936.90"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("E32").Value = "  -7.15%  "
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'8.91"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").Value = "'0.0991"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("E38").Value = "  +4.98%  "
$ws.Range("D39").Value = "'5.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'0.977"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("E41").Value = "  -4.90%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'43.84"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").Value = "'47.11"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").Value = "'151.98"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").Value = "'0.294"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").Value = "'8.28"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").Value = "'1.83"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "'26.66"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.91%  "
